$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to move from
# 45183 (2023-09-14) to 45184 (2023-09-15) for rows 2 through 11.
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45183) {
        $cell.Value2 = 45184
    }
}
